$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text in the source data (e.g.
# "22.139.30", "1.007"), which Excel would otherwise auto-convert to numbers.
# Force text storage for the edited range, write the values, then restore the
# default "Normal" style so no stray number-format survives on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.139.30'
$ws.Range("E2").Value = '  +7.78%  '
$ws.Range("D3").Value = '1.589.01'
$ws.Range("E3").Value = '  +7.91%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '0.9921'
$ws.Range("E5").Value = '  +3.59%  '
$ws.Range("D6").Value = '298.12'
$ws.Range("D7").Value = '0.3622'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '0.3340'
$ws.Range("E8").Value = '  +8.55%  '
$ws.Range("D9").Value = '40.93'
$ws.Range("D10").Value = '1.116'
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("D11").Value = '0.06944'
$ws.Range("E11").Value = '  +4.30%  '
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '19.37'
$ws.Range("E13").Value = '  +6.70%  '
$ws.Range("D14").Value = '5.820'
$ws.Range("E14").Value = '  +5.42%  '
$ws.Range("D15").Value = '6.528'
$ws.Range("E15").Value = '  +5.71%  '
$ws.Range("D16").Value = '0.9931'
$ws.Range("E16").Value = '  +3.63%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001060'
$ws.Range("E17").Value = '  +3.17%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '1.584.41'
$ws.Range("E18").Value = '  +7.56%  '
$ws.Range("D19").Value = '0.06576'
$ws.Range("E19").Value = '  +10.73%  '
$ws.Range("D20").Value = '75.99'
$ws.Range("E20").Value = '  +10.39%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '15.82'
$ws.Range("E21").Value = '  +8.75%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.914'
$ws.Range("E22").Value = '  +7.73%  '
$ws.Range("D23").Value = '11.61'
$ws.Range("E23").Value = '  +4.06%  '
$ws.Range("D24").Value = '22.175.04'
$ws.Range("E24").Value = '  +7.94%  '
$ws.Range("D25").Value = '2.376'
$ws.Range("E25").Value = '  +4.84%  '
$ws.Range("D26").Value = '2.496'
$ws.Range("E26").Value = '  +17.05%  '
$ws.Range("D27").Value = '148.05'
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("D28").Value = '19.14'
$ws.Range("E28").Value = '  +11.67%  '
$ws.Range("D29").Value = '1.754.76'
$ws.Range("E29").Value = '  +7.44%  '
$ws.Range("D30").Value = '121.77'
$ws.Range("E30").Value = '  +6.97%  '
$ws.Range("D31").Value = '3.980'
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("D32").Value = '5.887'
$ws.Range("E32").Value = '  +18.55%  '
$ws.Range("D33").Value = '0.9177'
$ws.Range("E33").Value = '  +14.04%  '
$ws.Range("D34").Value = '0.08131'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("D35").Value = '1.612'
$ws.Range("E35").Value = '  +6.31%  '
$ws.Range("D36").Value = '11.66'
$ws.Range("E36").Value = '  +12.17%  '
$ws.Range("D37").Value = '5.113'
$ws.Range("E37").Value = '  +8.06%  '
$ws.Range("D38").Value = '1.234'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").Value = '8.321'
$ws.Range("E39").Value = '  +12.32%  '
$ws.Range("D40").Value = '0.05976'
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("D41").Value = '0.02173'
$ws.Range("E41").Value = '  +5.74%  '
$ws.Range("D42").Value = '0.1981'
$ws.Range("E42").Value = '  +5.77%  '
$ws.Range("D43").Value = '0.9941'
$ws.Range("E43").Value = '  +3.71%  '
$ws.Range("D44").Value = '0.5778'
$ws.Range("E44").Value = '  +9.57%  '
$ws.Range("D45").Value = '3.764'
$ws.Range("E45").Value = '  +6.96%  '
$ws.Range("D46").Value = '12.85'
$ws.Range("E46").Value = '  +5.19%  '
$ws.Range("D47").Value = '125.03'
$ws.Range("E47").Value = '  +5.55%  '
$ws.Range("D48").Value = '0.5550'
$ws.Range("E48").Value = '  +6.62%  '
$ws.Range("D49").Value = '1.935'
$ws.Range("E49").Value = '  +6.53%  '
$ws.Range("D50").Value = '0.06706'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("D51").Value = '72.37'

$ws.Range("D2:D51").Style = "Normal"
